# Change the "maxDistancePreference" column's datatype from
# varchar(255) to float in the `patients` table DDL snippet.
#
# "maxDistancePreference " is a unique anchor in the document content,
# so locate it first, then replace the "varchar(255) not null," text
# that immediately follows it with "float        not null,".

$d = $word.ActiveDocument

$anchor = $d.Content.Duplicate
$found = $anchor.Find.Execute("maxDistancePreference ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find anchor text 'maxDistancePreference '"
}

$oldType = "varchar(255) not null,"
$target = $d.Range($anchor.End, $anchor.End + $oldType.Length)

if ($target.Text -ne $oldType) {
    throw "Unexpected text following anchor: [$($target.Text)]"
}

$target.Text = "float        not null,"
